$wb = $excel.ActiveWorkbook

# --- Add a new first sheet "EnvControl", keep the original sheet as "LoginData" ---
$loginData = $wb.Worksheets.Item(1)
$envControl = $wb.Worksheets.Add($loginData)
$envControl.Name = "EnvControl"

$loginData = $wb.Worksheets.Item("Sheet1")
$loginData.Name = "LoginData"

# --- Populate EnvControl with the environment execution-control table ---
$envControl.Range("A1").Value = "env"
$envControl.Range("B1").Value = "execute"
$envControl.Range("A2").Value = "uat"
$envControl.Range("B2").Value = "yes"
$envControl.Range("A3").Value = "qa"
$envControl.Range("B3").Value = "no"
$envControl.Range("A4").Value = "dev"
$envControl.Range("B4").Value = "no"

# Match the header style used on LoginData's header row (bold + yellow fill + centered)
$loginData.Range("A1:B1").Copy()
$envControl.Range("A1:B1").PasteSpecial(-4122)
$envControl.Application.CutCopyMode = $false

# --- LoginData data tweak: row 2's env list narrows from "qa,uat" to just "uat" ---
$loginData.Range("E2").Value = "uat"

# --- Selections / active tab to match the saved workbook state ---
[void]$envControl.Range("B4").Select()
[void]$loginData.Range("D2").Select()
[void]$loginData.Activate()
